# Updated symbol list on Sat Dec 24 05:45:05 UTC 2022 with GitHub Actions
#
# All values in this sheet are stored as text (inlineStr), even the
# numeric-looking price column D. Plain `.Value = "123.45"` assignment
# would be auto-coerced to a number by Excel, so every write goes through
# Set-Text which forces Text number-formatting for the duration of the
# write and then restores the cell to the default "Normal" style
# afterwards (avoids leaving a stray NumberFormat/style behind).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($a1, $value) {
    $cell = $ws.Range($a1)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- price-only updates (rows 2-17) ---
Set-Text "D2"  "245.24"
Set-Text "D3"  "22.04"
Set-Text "D4"  "5.332"
Set-Text "D5"  "0.05970"
Set-Text "D6"  "3.397"
Set-Text "D7"  "6.392"
Set-Text "D8"  "0.8124"
Set-Text "D9"  "0.9641"
Set-Text "D10" "0.1427"
Set-Text "D11" "0.07403"
Set-Text "D12" "0.03496"
Set-Text "D13" "0.03068"
Set-Text "D14" "0.09400"
Set-Text "D16" "0.001588"
Set-Text "D17" "0.04810"

# --- rows 18-24: coin ranking reshuffled up by one, with refreshed prices ---
Set-Text "B18" "TigerCash"
Set-Text "C18" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-Text "D18" "0.006234"
Set-Text "E18" "17TigerCashTCH"

Set-Text "B19" "HotbitToken"
Set-Text "C19" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-Text "D19" "0.004138"
Set-Text "E19" "18HotbitTokenHTB"

Set-Text "B20" "BitKan"
Set-Text "C20" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-Text "D20" "0.0009848"
Set-Text "E20" "19BitKanKAN"

Set-Text "B21" "NitroEx"
Set-Text "C21" "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-Text "D21" "0.00009708"
Set-Text "E21" "20NitroExNTX"

Set-Text "B22" "LEO"
Set-Text "C22" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-Text "D22" "3.747"
Set-Text "E22" "21LEOLEO"

Set-Text "B23" "BTSEToken"
Set-Text "C23" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-Text "D23" "2.166"
Set-Text "E23" "22BTSETokenBTSE"

Set-Text "B24" "One"
Set-Text "C24" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-Text "D24" "0.01120"
Set-Text "E24" "23OneONEBestin24h"

# --- price-only updates (rows 40-48) ---
Set-Text "D40" "0.03910"
Set-Text "D41" "0.006515"
Set-Text "D42" "0.1071"
Set-Text "D43" "0.003002"
Set-Text "D44" "0.005367"
Set-Text "D45" "0.00005316"

# row 47 lost its "Bestin24h" badge
Set-Text "E47" "46CoinbaseStockTokenCOIN"

# row 48 price update, and it gained the "Worstin24h" badge
Set-Text "D48" "0.03931"
Set-Text "E48" "47BOLOBOLOWorstin24h"
